$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Bobby Portis) is removed entirely; rows below shift up.
$ws.Rows(11).Delete()

# A new player, Josh Green, is inserted right before Giannis Antetokounmpo
# (which, after the deletion above, sits at row 16).
$ws.Rows(16).Insert()
$ws.Range("A16").Value = "Josh Green"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Charlotte Hornets"
